$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 101, shifting rows 101:212 down to 102:213
$ws.Rows("101:101").Insert()

# Populate the new row 101 with the new weekly record
$ws.Range("A101").Value = 8
$ws.Range("B101").Value = "Terminal La Palmera de La Serena"
$ws.Range("C101").Value = "Coquimbo"
$ws.Range("D101").Value = 44494
$ws.Range("E101").Value = 4
$ws.Range("F101").Value = 100114013
$ws.Range("G101").Value = "Zanahoria"
$ws.Range("H101").Value = "Sin especificar"
$ws.Range("I101").Value = "Primera"
$ws.Range("J101").Value = 600
$ws.Range("K101").Value = 6800
$ws.Range("L101").Value = 7000
$ws.Range("M101").Value = 6900
$ws.Range("N101").Value = "$/saco 20 kilos"
$ws.Range("O101").Value = "Provincia del Elquí"
$ws.Range("P101").Value = 345
$ws.Range("Q101").Value = 20
$ws.Range("R101").Value = "Hortaliza"

Write-Output "Done"
